$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove now-blank trailing rows 613:616 (their data moved into rows 586:612) ---

$ws.Range("A613:L616").EntireRow.Delete()

# --- Apply number formats by copying from existing template cells with identical styles ---
# A,B,C,D columns (styles 15,31,10,5) copied from row 578
$ws.Range("A578:D578").Copy()
$ws.Range("A586:D612").PasteSpecial(-4122)
# F column (style 16) copied from row 578
$ws.Range("F578").Copy()
$ws.Range("F586:F612").PasteSpecial(-4122)
# E column default style (5) copied from row 578 for all rows first ...
$ws.Range("E578").Copy()
$ws.Range("E586:E612").PasteSpecial(-4122)
# ... then the rows whose E uses style 14 get overwritten from row 572
$ws.Range("E572").Copy()
$ws.Range("E586").PasteSpecial(-4122)
$ws.Range("E592").PasteSpecial(-4122)
$ws.Range("E594").PasteSpecial(-4122)
$ws.Range("E597").PasteSpecial(-4122)
$ws.Range("E602").PasteSpecial(-4122)
$ws.Range("E606").PasteSpecial(-4122)
$ws.Range("E608").PasteSpecial(-4122)

# G,H,K columns (styles 13,18,5) copied from row 578 for the rows that carry summary formulas
$ws.Range("G578:H578").Copy()
$ws.Range("G586:H586").PasteSpecial(-4122)
$ws.Range("G589:H589").PasteSpecial(-4122)
$ws.Range("G591:H591").PasteSpecial(-4122)
$ws.Range("G594:H594").PasteSpecial(-4122)
$ws.Range("G597:H597").PasteSpecial(-4122)
$ws.Range("G600:H600").PasteSpecial(-4122)
$ws.Range("G602:H602").PasteSpecial(-4122)
$ws.Range("G605:H605").PasteSpecial(-4122)
$ws.Range("G608:H608").PasteSpecial(-4122)
$ws.Range("G611:H611").PasteSpecial(-4122)
$ws.Range("K578").Copy()
$ws.Range("K586").PasteSpecial(-4122)
$ws.Range("K589").PasteSpecial(-4122)
$ws.Range("K591").PasteSpecial(-4122)
$ws.Range("K594").PasteSpecial(-4122)
$ws.Range("K597").PasteSpecial(-4122)
$ws.Range("K600").PasteSpecial(-4122)
$ws.Range("K602").PasteSpecial(-4122)
$ws.Range("K605").PasteSpecial(-4122)
$ws.Range("K608").PasteSpecial(-4122)
$ws.Range("K611").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate values & formulas ---
# Row 586
$ws.Range("A586").Value2 = 'Jeremy'
$ws.Range("B586").Value2 = 45099
$ws.Range("C586").Value2 = 'Castello_south'
$ws.Range("D586").Value2 = 49.98
$ws.Range("E586").Value2 = 2715.04
$ws.Range("F586").Value2 = 'T2_t0_AMB_blank_01a'
$ws.Range("G586").Formula = '=AVERAGE(E587:E588)'
$ws.Range("H586").Formula = '=STDEV(E587:E588)'
$ws.Range("K586").Formula = '=G586+$J$575'

# Row 587
$ws.Range("A587").Value2 = 'Jeremy'
$ws.Range("B587").Value2 = 45099
$ws.Range("C587").Value2 = 'Castello_south'
$ws.Range("D587").Value2 = 50.43
$ws.Range("E587").Value2 = 2692.92
$ws.Range("F587").Value2 = 'T2_t0_AMB_blank_01b'

# Row 588
$ws.Range("A588").Value2 = 'Jeremy'
$ws.Range("B588").Value2 = 45099
$ws.Range("C588").Value2 = 'Castello_south'
$ws.Range("D588").Value2 = 50.15
$ws.Range("E588").Value2 = 2693.4
$ws.Range("F588").Value2 = 'T2_t0_AMB_blank_01c'

# Row 589
$ws.Range("A589").Value2 = 'Jeremy'
$ws.Range("B589").Value2 = 45099
$ws.Range("C589").Value2 = 'Castello_south'
$ws.Range("D589").Value2 = 50.55
$ws.Range("E589").Value2 = 2686.54
$ws.Range("F589").Value2 = 'T2_t0_AMB_blank_02a'
$ws.Range("G589").Formula = '=AVERAGE(E589:E590)'
$ws.Range("H589").Formula = '=STDEV(E589:E590)'
$ws.Range("K589").Formula = '=G589+$J$575'

# Row 590
$ws.Range("A590").Value2 = 'Jeremy'
$ws.Range("B590").Value2 = 45099
$ws.Range("C590").Value2 = 'Castello_south'
$ws.Range("D590").Value2 = 50.75
$ws.Range("E590").Value2 = 2681.8
$ws.Range("F590").Value2 = 'T2_t0_AMB_blank_02b'

# Row 591
$ws.Range("A591").Value2 = 'Jeremy'
$ws.Range("B591").Value2 = 45099
$ws.Range("C591").Value2 = 'Castello_south'
$ws.Range("D591").Value2 = 49.37
$ws.Range("E591").Value2 = 2538.7600000000002
$ws.Range("F591").Value2 = 'T2_t1_AMB_tile_01a'
$ws.Range("G591").Formula = '=AVERAGE(E591,E593)'
$ws.Range("H591").Formula = '=STDEV(E591,E593)'
$ws.Range("K591").Formula = '=G591+$J$575'

# Row 592
$ws.Range("A592").Value2 = 'Jeremy'
$ws.Range("B592").Value2 = 45099
$ws.Range("C592").Value2 = 'Castello_south'
$ws.Range("D592").Value2 = 50.21
$ws.Range("E592").Value2 = 2524.25
$ws.Range("F592").Value2 = 'T2_t1_AMB_tile_01b'

# Row 593
$ws.Range("A593").Value2 = 'Jeremy'
$ws.Range("B593").Value2 = 45099
$ws.Range("C593").Value2 = 'Castello_south'
$ws.Range("D593").Value2 = 48.43
$ws.Range("E593").Value2 = 2539.6799999999998
$ws.Range("F593").Value2 = 'T2_t1_AMB_tile_01c'

# Row 594
$ws.Range("A594").Value2 = 'Jeremy'
$ws.Range("B594").Value2 = 45099
$ws.Range("C594").Value2 = 'Castello_south'
$ws.Range("D594").Value2 = 50.41
$ws.Range("E594").Value2 = 2578.7800000000002
$ws.Range("F594").Value2 = 'T2_t1_AMB_tile_02a'
$ws.Range("G594").Formula = '=AVERAGE(E595:E596)'
$ws.Range("H594").Formula = '=STDEV(E595:E596)'
$ws.Range("K594").Formula = '=G594+$J$575'

# Row 595
$ws.Range("A595").Value2 = 'Jeremy'
$ws.Range("B595").Value2 = 45099
$ws.Range("C595").Value2 = 'Castello_south'
$ws.Range("D595").Value2 = 49.38
$ws.Range("E595").Value2 = 2592.5
$ws.Range("F595").Value2 = 'T2_t1_AMB_tile_02c'

# Row 596
$ws.Range("A596").Value2 = 'Jeremy'
$ws.Range("B596").Value2 = 45099
$ws.Range("C596").Value2 = 'Castello_south'
$ws.Range("D596").Value2 = 49.42
$ws.Range("E596").Value2 = 2588.16
$ws.Range("F596").Value2 = 'T2_t1_AMB_tile_02e'

# Row 597
$ws.Range("A597").Value2 = 'Jeremy'
$ws.Range("B597").Value2 = 45099
$ws.Range("C597").Value2 = 'Castello_south'
$ws.Range("D597").Value2 = 51.38
$ws.Range("E597").Value2 = 2545.4299999999998
$ws.Range("F597").Value2 = 'T2_t1_AMB_tile_03a'
$ws.Range("G597").Formula = '=AVERAGE(E598:E599)'
$ws.Range("H597").Formula = '=STDEV(E598:E599)'
$ws.Range("K597").Formula = '=G597+$J$575'

# Row 598
$ws.Range("A598").Value2 = 'Jeremy'
$ws.Range("B598").Value2 = 45099
$ws.Range("C598").Value2 = 'Castello_south'
$ws.Range("D598").Value2 = 51.28
$ws.Range("E598").Value2 = 2560.7199999999998
$ws.Range("F598").Value2 = 'T2_t1_AMB_tile_03b'

# Row 599
$ws.Range("A599").Value2 = 'Jeremy'
$ws.Range("B599").Value2 = 45099
$ws.Range("C599").Value2 = 'Castello_south'
$ws.Range("D599").Value2 = 50.12
$ws.Range("E599").Value2 = 2556.3200000000002
$ws.Range("F599").Value2 = 'T2_t1_AMB_tile_03c'

# Row 600
$ws.Range("A600").Value2 = 'Jeremy'
$ws.Range("B600").Value2 = 45099
$ws.Range("C600").Value2 = 'Castello_south'
$ws.Range("D600").Value2 = 51.33
$ws.Range("E600").Value2 = 2680.89
$ws.Range("F600").Value2 = 'T2_t1_AMB_blank_01a'
$ws.Range("G600").Formula = '=AVERAGE(E600:E601)'
$ws.Range("H600").Formula = '=STDEV(E600:E601)'
$ws.Range("K600").Formula = '=G600+$J$575'

# Row 601
$ws.Range("A601").Value2 = 'Jeremy'
$ws.Range("B601").Value2 = 45099
$ws.Range("C601").Value2 = 'Castello_south'
$ws.Range("D601").Value2 = 51.22
$ws.Range("E601").Value2 = 2676.85
$ws.Range("F601").Value2 = 'T2_t1_AMB_blank_01b'

# Row 602
$ws.Range("A602").Value2 = 'Jeremy'
$ws.Range("B602").Value2 = 45099
$ws.Range("C602").Value2 = 'Castello_south'
$ws.Range("D602").Value2 = 53.6
$ws.Range("E602").Value2 = 2626.42
$ws.Range("F602").Value2 = 'T2_t1_AMB_tile_04a'
$ws.Range("G602").Formula = '=AVERAGE(E603:E604)'
$ws.Range("H602").Formula = '=STDEV(E603:E604)'
$ws.Range("K602").Formula = '=G602+$J$575'

# Row 603
$ws.Range("A603").Value2 = 'Jeremy'
$ws.Range("B603").Value2 = 45099
$ws.Range("C603").Value2 = 'Castello_south'
$ws.Range("D603").Value2 = 51.26
$ws.Range("E603").Value2 = 2613.58
$ws.Range("F603").Value2 = 'T2_t1_AMB_tile_04b'

# Row 604
$ws.Range("A604").Value2 = 'Jeremy'
$ws.Range("B604").Value2 = 45099
$ws.Range("C604").Value2 = 'Castello_south'
$ws.Range("D604").Value2 = 50.69
$ws.Range("E604").Value2 = 2619.02
$ws.Range("F604").Value2 = 'T2_t1_AMB_tile_04c'

# Row 605
$ws.Range("A605").Value2 = 'Jeremy'
$ws.Range("B605").Value2 = 45099
$ws.Range("C605").Value2 = 'Castello_south'
$ws.Range("D605").Value2 = 53.95
$ws.Range("E605").Value2 = 2543.17
$ws.Range("F605").Value2 = 'T2_t1_AMB_tile_05a'
$ws.Range("G605").Formula = '=AVERAGE(E605,E607)'
$ws.Range("H605").Formula = '=STDEV(E605,E607)'
$ws.Range("K605").Formula = '=G605+$J$575'

# Row 606
$ws.Range("A606").Value2 = 'Jeremy'
$ws.Range("B606").Value2 = 45099
$ws.Range("C606").Value2 = 'Castello_south'
$ws.Range("D606").Value2 = 51
$ws.Range("E606").Value2 = 2532.25
$ws.Range("F606").Value2 = 'T2_t1_AMB_tile_05b'

# Row 607
$ws.Range("A607").Value2 = 'Jeremy'
$ws.Range("B607").Value2 = 45099
$ws.Range("C607").Value2 = 'Castello_south'
$ws.Range("D607").Value2 = 50.84
$ws.Range("E607").Value2 = 2548.73
$ws.Range("F607").Value2 = 'T2_t1_AMB_tile_05c'

# Row 608
$ws.Range("A608").Value2 = 'Jeremy'
$ws.Range("B608").Value2 = 45099
$ws.Range("C608").Value2 = 'Castello_south'
$ws.Range("D608").Value2 = 50.39
$ws.Range("E608").Value2 = 2552.04
$ws.Range("F608").Value2 = 'T2_t1_AMB_tile_06a'
$ws.Range("G608").Formula = '=AVERAGE(E609:E610)'
$ws.Range("H608").Formula = '=STDEV(E609:E610)'
$ws.Range("K608").Formula = '=G608+$J$575'

# Row 609
$ws.Range("A609").Value2 = 'Jeremy'
$ws.Range("B609").Value2 = 45099
$ws.Range("C609").Value2 = 'Castello_south'
$ws.Range("D609").Value2 = 49.68
$ws.Range("E609").Value2 = 2593.4699999999998
$ws.Range("F609").Value2 = 'T2_t1_AMB_tile_06b'

# Row 610
$ws.Range("A610").Value2 = 'Jeremy'
$ws.Range("B610").Value2 = 45099
$ws.Range("C610").Value2 = 'Castello_south'
$ws.Range("D610").Value2 = 50.42
$ws.Range("E610").Value2 = 2596.79
$ws.Range("F610").Value2 = 'T2_t1_AMB_tile_06e'

# Row 611
$ws.Range("A611").Value2 = 'Jeremy'
$ws.Range("B611").Value2 = 45099
$ws.Range("C611").Value2 = 'Castello_south'
$ws.Range("D611").Value2 = 51.14
$ws.Range("E611").Value2 = 2689.54
$ws.Range("F611").Value2 = 'T2_t1_AMB_blank_02a'
$ws.Range("G611").Formula = '=AVERAGE(E611:E612)'
$ws.Range("H611").Formula = '=STDEV(E611:E612)'
$ws.Range("K611").Formula = '=G611+$J$575'

# Row 612
$ws.Range("A612").Value2 = 'Jeremy'
$ws.Range("B612").Value2 = 45099
$ws.Range("C612").Value2 = 'Castello_south'
$ws.Range("D612").Value2 = 50.5
$ws.Range("E612").Value2 = 2685.28
$ws.Range("F612").Value2 = 'T2_t1_AMB_blank_02b'

# --- Update sheet view selection/scroll to match final state ---
$ws.Range("G596").Select()
